# The commit swaps the theme used as the deck's main design ("Integral")
# with the plain default "Office Theme" that previously only backed the
# notes master: theme1.xml (the SlideMaster/Design theme that PowerPoint's
# object model actually exposes for editing) ends up holding the "Office
# Theme" palette, while the content that used to live in theme1.xml
# ("Integral") moves to theme2.xml.
#
# The fontScheme and fmtScheme blocks are byte-for-byte identical between
# the two theme parts already, and dk1/lt1 (black/white) are identical
# too, so the only real work is recoloring the ten clrScheme entries that
# differ (dk2, lt2, accent1-6, hlink, folHlink) on the editable theme to
# the "Office Theme" values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Set-SchemeColor {
    param($themeColors, $index, $hex)
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $themeColors.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Index map (ThemeColorScheme is 1-based, 12 entries):
#  1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#  8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
Set-SchemeColor $tcs 3  "44546A"   # dk2
Set-SchemeColor $tcs 4  "E7E6E6"   # lt2
Set-SchemeColor $tcs 5  "5B9BD5"   # accent1
Set-SchemeColor $tcs 6  "ED7D31"   # accent2
Set-SchemeColor $tcs 7  "A5A5A5"   # accent3
Set-SchemeColor $tcs 8  "FFC000"   # accent4
Set-SchemeColor $tcs 9  "4472C4"   # accent5
Set-SchemeColor $tcs 10 "70AD47"   # accent6
Set-SchemeColor $tcs 11 "0563C1"   # hlink
Set-SchemeColor $tcs 12 "954F72"   # folHlink
